# Updated values with more steps
# - Renames Sheet1 -> VoltageRailSpec
# - Adds a new "Task4_Efficiency" worksheet after it, with the Task 4
#   efficiency-vs-load-resistance data table (R_load, V_in, I_in, V_out,
#   I_out, Efficiency) and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Rename the original sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "VoltageRailSpec"

# Insert the new sheet right after the renamed one.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Task4_Efficiency"

# Headers -- written in this particular column order so the shared-string
# table ends up populated in the same sequence as the source workbook.
$ws2.Range("C1").Value = '$I_{in}$ (A)'
$ws2.Range("D1").Value = '$V_{out}$ (V)'
$ws2.Range("E1").Value = '$I_{out}$ (A)'
$ws2.Range("F1").Value = 'Efficiency $\eta$'
$ws2.Range("B1").Value = '$V_{in}$ (V)'
$ws2.Range("A1").Value = '$R_{load}$ ($\Omega$)'

# Data rows: A=R_load, B=V_in, C=I_in, D=V_out, E=I_out, F=Efficiency (=D*E/C/B)
$data = @(
    @(2.47,   100, 92.0614,  88.3454, 35.7674),
    @(12.35,  100, 61.2578,  179.641, 14.5458),
    @(24.7,   100, 44.8744,  259.556, 10.5083),
    @(123.5,  100, 13.1896,  377.67,   3.05805),
    @(148.2,  100, 11.3059,  388.214,  2.61953),
    @(172.9,  100, 9.83076,  393.375,  2.27516),
    @(197.6,  100, 8.69834,  397.339,  2.01082),
    @(222.3,  100, 7.80164,  400.479,  1.80152),
    @(247,    100, 7.074,    403.027,  1.63169)
)

$row = 2
foreach ($r in $data) {
    $ws2.Range("A$row").Value = $r[0]
    $ws2.Range("B$row").Value = $r[1]
    $ws2.Range("C$row").Value = $r[2]
    $ws2.Range("D$row").Value = $r[3]
    $ws2.Range("E$row").Value = $r[4]
    $ws2.Range("F$row").Formula = "=D$row*E$row/C$row/B$row"
    $row++
}

# Column sizing to roughly match the author's best-fit widths.
$ws2.Columns.Item(1).ColumnWidth = 19.7
$ws2.Columns.Item(2).ColumnWidth = 9.8
$ws2.Columns.Item(3).ColumnWidth = 9.3
$ws2.Columns.Item(5).ColumnWidth = 10.6
$ws2.Columns.Item(6).ColumnWidth = 14.1

# View settings: selection, zoom and making the new sheet the active tab.
$ws2.Range("A1:F10").Select()
$ws2.Activate()
$win = $wb.Windows.Item(1)
$win.Zoom = 140
